$wb = $excel.ActiveWorkbook

# --- Update the driver test data (email / mobile) on "driver_Post" ---
$wsPost = $wb.Worksheets.Item("driver_Post")
$wsPost.Range("B4").Value = "driver.test_51@gmail.com"
$wsPost.Range("B5").Value = "7755663251"

# --- Update the same driver test data on "driver_Post_ValidationData" ---
$wsVal = $wb.Worksheets.Item("driver_Post_ValidationData")
$wsVal.Range("B5").Value = "driver.test_51@gmail.com"
$wsVal.Range("B6").Value = "7755663251"

# --- Update the same driver test data on "driver_verify_Post" ---
$wsVerify = $wb.Worksheets.Item("driver_verify_Post")
$wsVerify.Range("B2").Value = "driver.test_51@gmail.com"

# --- Move the active/selected tab from "driver_Post" to "driver_Post_ValidationData" ---
$wsVal.Activate()
$wsVal.Range("B5").Select()
